# Add a new "Desc" column (G) to the parts table with a short marketing
# description for every part, in the same order the original author typed
# them in (a handful of CPU-cooler / motherboard rows were filled in last).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 'Desc'
$ws.Range("G3").Value = 'Large chassi with two large fans on the front and one on the back'
$ws.Range("G4").Value = 'Large chassi with three fans on the front'
$ws.Range("G5").Value = 'Dark large chassi with three medium-sized fans on the front'
$ws.Range("G6").Value = 'Nice transparent chassis with great potential for nice computers'
$ws.Range("G7").Value = 'Popular midtower chassis with nice design'
$ws.Range("G8").Value = 'Fresh midtower chassis'
$ws.Range("G10").Value = 'AMD''s entry processor, good for less demanding applications'
$ws.Range("G11").Value = 'AMD''s midrange processor, good cpu for the money'
$ws.Range("G12").Value = 'AMD''s higher range of processors, good for editing and gaming'
$ws.Range("G13").Value = 'Intel mid range processor, good for gaming'
$ws.Range("G14").Value = 'Intel high range processor, good for editing'
$ws.Range("G15").Value = 'Intel''s best processor, awesome for games and editing etc.'
$ws.Range("G17").Value = 'A standard and simple CPU fan that fixes the job'
$ws.Range("G18").Value = 'Nice fan from be quiet, which is also very quiet'
$ws.Range("G19").Value = 'Nice CPU fan for a good price'
$ws.Range("G22").Value = 'Simple but good water-cooled CPU fan'
$ws.Range("G24").Value = 'Good entrance card for games from AMD'
$ws.Range("G25").Value = 'Good graphics card for the money from AMD'
$ws.Range("G26").Value = 'AMD''s top-rated graphics card that crushes all games'
$ws.Range("G27").Value = 'Nvidia''s cheapest and best graphics card for the money'
$ws.Range("G28").Value = 'Nice graphics card from Nvidia, perfect for middle range computers'
$ws.Range("G29").Value = 'Awesome graphics card for those who are new to RTX'
$ws.Range("G30").Value = 'Nvidia''s best 20 series graphics card'
$ws.Range("G31").Value = 'A perfect graphics card that crushes all AAA games with relief'
$ws.Range("G32").Value = 'Nvidia''s best graphics card, which tremendously crushes all competitors'
$ws.Range("G34").Value = 'Nice motherboard for AMD socket processors'
$ws.Range("G36").Value = 'Nice motherboard from ASUS for Intel socket processors'
$ws.Range("G39").Value = 'Standard version of Windows, supports up to 128GB RAM'
$ws.Range("G40").Value = 'Top version of Windows, supports up to 2TB of RAM'
$ws.Range("G42").Value = '500W power supply from EVGA'
$ws.Range("G43").Value = 'Reliable 600W power supply from EVGA'
$ws.Range("G44").Value = '750W power supply that can handle everything from Corsair'
$ws.Range("G45").Value = '850W power supply from Corsair, you will never need more watts'
$ws.Range("G47").Value = '2x4 - 8GB DDR4-3000Mhz RAM from Corsair'
$ws.Range("G48").Value = '2x8 - 16GB DDR4-3200Mhz RAM from Corsair'
$ws.Range("G49").Value = '2x16 - 32GB DDR4-3200Mhz RAM from Corsair'
$ws.Range("G50").Value = '4x8 - 32GB DDR4-3600Mhz RAM from Corsair'
$ws.Range("G51").Value = '2x16 - 32GB DDR4-3600Mhz RAM from G.Skill'
$ws.Range("G52").Value = '2x8 - 16GB DDR4-3600Mhz RAM from G.Skill'
$ws.Range("G54").Value = 'Fresh and small 240GB SSD from Kingston'
$ws.Range("G55").Value = 'Nice 1TB SSD from Samsung, perfect second hard drive for your computer'
$ws.Range("G56").Value = 'Good SSD of 500GB from Samsung'
$ws.Range("G57").Value = '1TB M.2-NVME SSD from Samsung with insanely fast speeds, perfect for your gaming computer'
$ws.Range("G58").Value = 'A 500GB M.2-NVME SSD from Samsung that works perfectly for all computers'
$ws.Range("G59").Value = 'Mechanical hard drive from Seagate at 1TB, perfect for editing savings'
$ws.Range("G60").Value = 'A mechanical hard drive with a large space of 2TB'
$ws.Range("G20").Value = 'Good water-cooled CPU fan, even with nice rgb'
$ws.Range("G21").Value = 'An awesome CPU fan powered on water cooling with tripple fans for efficient cooling'
$ws.Range("G35").Value = 'Perfect motherboard from ASUS for Intel socket processors'
$ws.Range("G37").Value = 'Nice motherboard for AMD socket processors'

# Set column G width to match the authored layout (closest value achievable
# given the engine's 1/6-character width snapping).
$ws.Columns.Item(7).ColumnWidth = 85.43

# Restore the final selection/scroll state from the saved workbook.
$ws.Range("G35").Select()
